$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.736.44"
$ws.Range("D2").Style = $ws.Range("D4").Style
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "'1.881.09"
$ws.Range("D3").Style = $ws.Range("D4").Style
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "'333.16"
$ws.Range("D5").Style = $ws.Range("D4").Style
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.4709"
$ws.Range("D7").Style = $ws.Range("D4").Style
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("D8").Value = "'0.3935"
$ws.Range("D8").Style = $ws.Range("D4").Style
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").Value = "'47.52"
$ws.Range("D9").Style = $ws.Range("D4").Style
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "'0.08071"
$ws.Range("D10").Style = $ws.Range("D4").Style
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").Value = "'1.029"
$ws.Range("D11").Style = $ws.Range("D4").Style
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "'22.22"
$ws.Range("D12").Style = $ws.Range("D4").Style
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("D13").Value = "'1.886.07"
$ws.Range("D13").Style = $ws.Range("D4").Style
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "'5.981"
$ws.Range("D14").Style = $ws.Range("D4").Style
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "'7.144"
$ws.Range("D15").Style = $ws.Range("D4").Style
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "'1.010"
$ws.Range("D16").Style = $ws.Range("D4").Style
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.00001049"
$ws.Range("D17").Style = $ws.Range("D4").Style
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.06704"
$ws.Range("D18").Style = $ws.Range("D4").Style
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'87.11"
$ws.Range("D19").Style = $ws.Range("D4").Style
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'17.35"
$ws.Range("D20").Style = $ws.Range("D4").Style
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.545"
$ws.Range("D22").Style = $ws.Range("D4").Style
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "'27.726.82"
$ws.Range("D23").Style = $ws.Range("D4").Style
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "'11.04"
$ws.Range("D24").Style = $ws.Range("D4").Style
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").Value = "'2.108.34"
$ws.Range("D26").Style = $ws.Range("D4").Style
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("D27").Value = "'160.04"
$ws.Range("D27").Style = $ws.Range("D4").Style
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("D28").Value = "'20.18"
$ws.Range("D28").Style = $ws.Range("D4").Style
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "'2.107"
$ws.Range("D29").Style = $ws.Range("D4").Style
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").Value = "'5.592"
$ws.Range("D30").Style = $ws.Range("D4").Style
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("D31").Value = "'121.94"
$ws.Range("D31").Style = $ws.Range("D4").Style
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "'0.9856"
$ws.Range("D32").Style = $ws.Range("D4").Style
$ws.Range("E32").Value = "  +3.94%  "
$ws.Range("D33").Value = "'0.09483"
$ws.Range("D33").Style = $ws.Range("D4").Style
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").Value = "'1.453"
$ws.Range("D34").Style = $ws.Range("D4").Style
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Value = "'3.623"
$ws.Range("D35").Style = $ws.Range("D4").Style
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").Value = "'5.362"
$ws.Range("D36").Style = $ws.Range("D4").Style
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").Value = "'0.06139"
$ws.Range("D37").Style = $ws.Range("D4").Style
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("D38").Value = "'0.02266"
$ws.Range("D38").Style = $ws.Range("D4").Style
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").Value = "'1.233"
$ws.Range("D39").Style = $ws.Range("D4").Style
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").Value = "'8.141"
$ws.Range("D40").Style = $ws.Range("D4").Style
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").Value = "'0.6007"
$ws.Range("D41").Style = $ws.Range("D4").Style
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D43").Value = "'10.29"
$ws.Range("D43").Style = $ws.Range("D4").Style
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("D44").Value = "'0.5733"
$ws.Range("D44").Style = $ws.Range("D4").Style
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").Value = "'12.26"
$ws.Range("D46").Style = $ws.Range("D4").Style
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("D47").Value = "'1.948"
$ws.Range("D47").Style = $ws.Range("D4").Style
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").Value = "'3.399"
$ws.Range("D48").Style = $ws.Range("D4").Style
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").Value = "'0.06903"
$ws.Range("D49").Style = $ws.Range("D4").Style
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("E50").Value = "  +5.70%  "
$ws.Range("E51").Value = "  +1.73%  "
